$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")
$ws.Copy($null, $ws)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Page2"
Write-Host $wb.Worksheets.Count
foreach ($s in $wb.Worksheets) { Write-Host $s.Name }

# Test row insert
$ws.Rows("12:13").Insert()
Write-Host "Inserted rows"
